$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1989.875
$ws.Range("I100").Value = 2067.5
$ws.Range("J100").Value = 1912.25
$ws.Range("K100").Value = 2067.5
$ws.Range("L100").Value = 1912.25
$ws.Range("M100").Value = -1526.5
$ws.Range("N100").Value = -2994.25
$ws.Range("H107").Value = 677.2
$ws.Range("I107").Value = 692.9167
$ws.Range("J107").Value = 300
$ws.Range("K107").Value = 692.9167
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1227.0833
$ws.Range("N107").Value = -4140
$ws.Range("H120").Value = 48251.25
$ws.Range("J120").Value = 48251.25
$ws.Range("L120").Value = 48251.25
$ws.Range("N120").Value = -57927.25
$ws.Range("H123").Value = 31353.867
$ws.Range("J123").Value = 31353.867
$ws.Range("L123").Value = 31353.867
$ws.Range("N123").Value = -41153.867
$ws.Range("H137").Value = 2766.899
$ws.Range("I137").Value = 1021.4167
$ws.Range("J137").Value = 3411.3845
$ws.Range("K137").Value = 3064.2501
$ws.Range("L137").Value = 10234.1535
$ws.Range("M137").Value = -514.2501000000002
$ws.Range("N137").Value = -15334.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 4500
$ws.Range("I36").Value = 4500
$ws.Range("K36").Value = 4500
$ws.Range("M36").Value = -4154
$ws.Range("H61").Value = 3577.3
$ws.Range("I61").Value = 3626.6924
$ws.Range("J61").Value = 3485.5715
$ws.Range("K61").Value = 3626.6924
$ws.Range("L61").Value = 3485.5715
$ws.Range("M61").Value = -3414.6924
$ws.Range("N61").Value = -3909.5715
$ws.Range("H74").Value = 1760.591
$ws.Range("I74").Value = 1078.6666
$ws.Range("J74").Value = 2578.9
$ws.Range("K74").Value = 1078.6666
$ws.Range("L74").Value = 2578.9
$ws.Range("M74").Value = -204.6666
$ws.Range("N74").Value = -4326.9
$ws.Range("H77").Value = 1760.591
$ws.Range("I77").Value = 1078.6666
$ws.Range("J77").Value = 2578.9
$ws.Range("K77").Value = 5393.333000000001
$ws.Range("L77").Value = 12894.5
$ws.Range("M77").Value = -1025.333000000001
$ws.Range("N77").Value = -21630.5
$ws.Range("H107").Value = 38000
$ws.Range("J107").Value = 38000
$ws.Range("L107").Value = 38000
$ws.Range("N107").Value = -45680
$ws.Range("H118").Value = 46666.668
$ws.Range("J118").Value = 46666.668
$ws.Range("L118").Value = 46666.668
$ws.Range("N118").Value = -49980.668
$ws.Range("H122").Value = 3576
$ws.Range("I122").Value = 3900
$ws.Range("K122").Value = 11700
$ws.Range("M122").Value = -9250
$ws.Range("H136").Value = 3577.3
$ws.Range("I136").Value = 3626.6924
$ws.Range("J136").Value = 3485.5715
$ws.Range("K136").Value = 10880.0772
$ws.Range("L136").Value = 10456.7145
$ws.Range("M136").Value = -8330.0772
$ws.Range("N136").Value = -15556.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 819.8182
$ws.Range("I94").Value = 777.25
$ws.Range("J94").Value = 933.3333
$ws.Range("K94").Value = 777.25
$ws.Range("L94").Value = 933.3333
$ws.Range("M94").Value = -326.25
$ws.Range("N94").Value = -1835.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 43745.168
$ws.Range("J20").Value = 43745.168
$ws.Range("L20").Value = 43745.168
$ws.Range("N20").Value = -44217.168
$ws.Range("H30").Value = 43745.168
$ws.Range("J30").Value = 43745.168
$ws.Range("L30").Value = 43745.168
$ws.Range("N30").Value = -43927.168
$ws.Range("H31").Value = 6559.048
$ws.Range("I31").Value = 3653.111
$ws.Range("J31").Value = 7351.5757
$ws.Range("K31").Value = 3653.111
$ws.Range("L31").Value = 7351.5757
$ws.Range("M31").Value = -3358.111
$ws.Range("N31").Value = -7941.5757
$ws.Range("H34").Value = 6559.048
$ws.Range("I34").Value = 3653.111
$ws.Range("J34").Value = 7351.5757
$ws.Range("K34").Value = 3653.111
$ws.Range("L34").Value = 7351.5757
$ws.Range("M34").Value = -3451.111
$ws.Range("N34").Value = -7755.5757
$ws.Range("H58").Value = 2109.4822
$ws.Range("I58").Value = 1886.1777
$ws.Range("J58").Value = 3023
$ws.Range("K58").Value = 1886.1777
$ws.Range("L58").Value = 3023
$ws.Range("M58").Value = -1683.1777
$ws.Range("N58").Value = -3429
$ws.Range("H99").Value = 2057.1035
$ws.Range("I99").Value = 1979.7142
$ws.Range("J99").Value = 2129.3333
$ws.Range("K99").Value = 1979.7142
$ws.Range("L99").Value = 2129.3333
$ws.Range("M99").Value = -481.7141999999999
$ws.Range("N99").Value = -5125.3333
$ws.Range("H126").Value = 2057.1035
$ws.Range("I126").Value = 1979.7142
$ws.Range("J126").Value = 2129.3333
$ws.Range("K126").Value = 5939.142599999999
$ws.Range("L126").Value = 6387.999899999999
$ws.Range("M126").Value = -3469.142599999999
$ws.Range("N126").Value = -11327.9999
$ws.Range("H128").Value = 43745.168
$ws.Range("J128").Value = 43745.168
$ws.Range("L128").Value = 43745.168
$ws.Range("N128").Value = -53705.168
$ws.Range("H132").Value = 28722.23
$ws.Range("I132").Value = 1361.8605
$ws.Range("J132").Value = 159444
$ws.Range("K132").Value = 4085.5815
$ws.Range("L132").Value = 478332
$ws.Range("M132").Value = -1555.5815
$ws.Range("N132").Value = -483392
$ws.Range("H134").Value = 3651.3572
$ws.Range("I134").Value = 1731.4286
$ws.Range("J134").Value = 5571.2856
$ws.Range("K134").Value = 5194.2858
$ws.Range("L134").Value = 16713.8568
$ws.Range("M134").Value = -2659.2858
$ws.Range("N134").Value = -21783.8568
$ws.Range("H136").Value = 2109.4822
$ws.Range("I136").Value = 1886.1777
$ws.Range("J136").Value = 3023
$ws.Range("K136").Value = 5658.5331
$ws.Range("L136").Value = 9069
$ws.Range("M136").Value = -3108.5331
$ws.Range("N136").Value = -14169

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 8607466
$ws.Range("I2").Value = 43.066666
$ws.Range("J2").Value = 17214888
$ws.Range("K2").Value = 258.399996
$ws.Range("L2").Value = 103289328
$ws.Range("M2").Value = -145.399996
$ws.Range("N2").Value = -103289554
$ws.Range("H4").Value = 2309.7097
$ws.Range("I4").Value = 88.42856999999999
$ws.Range("J4").Value = 2957.5833
$ws.Range("K4").Value = 265.28571
$ws.Range("L4").Value = 8872.749899999999
$ws.Range("M4").Value = -153.28571
$ws.Range("N4").Value = -9096.749899999999
$ws.Range("H64").Value = 3824.9375
$ws.Range("I64").Value = 2425
$ws.Range("J64").Value = 4291.5835
$ws.Range("K64").Value = 7275
$ws.Range("L64").Value = 12874.7505
$ws.Range("M64").Value = -7005
$ws.Range("N64").Value = -13414.7505
$ws.Range("H67").Value = 3824.9375
$ws.Range("I67").Value = 2425
$ws.Range("J67").Value = 4291.5835
$ws.Range("K67").Value = 7275
$ws.Range("L67").Value = 12874.7505
$ws.Range("M67").Value = -6339
$ws.Range("N67").Value = -14746.7505
$ws.Range("H70").Value = 4945.6665
$ws.Range("I70").Value = 3006
$ws.Range("J70").Value = 5499.857
$ws.Range("K70").Value = 9018
$ws.Range("L70").Value = 16499.571
$ws.Range("M70").Value = -8703
$ws.Range("N70").Value = -17129.571
$ws.Range("H73").Value = 4945.6665
$ws.Range("I73").Value = 3006
$ws.Range("J73").Value = 5499.857
$ws.Range("K73").Value = 9018
$ws.Range("L73").Value = 16499.571
$ws.Range("M73").Value = -7926
$ws.Range("N73").Value = -18683.571
$ws.Range("H75").Value = 2312.5557
$ws.Range("I75").Value = 956.5
$ws.Range("J75").Value = 2700
$ws.Range("K75").Value = 2869.5
$ws.Range("L75").Value = 8100
$ws.Range("M75").Value = -1871.5
$ws.Range("N75").Value = -10096
$ws.Range("H78").Value = 2312.5557
$ws.Range("I78").Value = 956.5
$ws.Range("J78").Value = 2700
$ws.Range("K78").Value = 8608.5
$ws.Range("L78").Value = 24300
$ws.Range("M78").Value = -3616.5
$ws.Range("N78").Value = -34284
$ws.Range("H117").Value = 1709.75
$ws.Range("I117").Value = 1333.3334
$ws.Range("J117").Value = 1935.6
$ws.Range("K117").Value = 4000.0002
$ws.Range("L117").Value = 5806.799999999999
$ws.Range("M117").Value = -558.0001999999999
$ws.Range("N117").Value = -12690.8
$ws.Range("H131").Value = 64896.21
$ws.Range("I131").Value = 20471.8
$ws.Range("J131").Value = 72829.14
$ws.Range("K131").Value = 61415.39999999999
$ws.Range("L131").Value = 218487.42
$ws.Range("M131").Value = -56375.39999999999
$ws.Range("N131").Value = -228567.42

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 11991.182
$ws.Range("J92").Value = 11991.182
$ws.Range("L92").Value = 11991.182
$ws.Range("N92").Value = -15735.182
$ws.Range("H132").Value = 2144.8462
$ws.Range("I132").Value = 1815.1389
$ws.Range("J132").Value = 2886.6875
$ws.Range("K132").Value = 5445.4167
$ws.Range("L132").Value = 8660.0625
$ws.Range("M132").Value = -2915.4167
$ws.Range("N132").Value = -13720.0625
$ws.Range("H136").Value = 23961.375
$ws.Range("J136").Value = 23961.375
$ws.Range("L136").Value = 71884.125
$ws.Range("N136").Value = -76984.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H111").Value = 44379
$ws.Range("J111").Value = 44379
$ws.Range("L111").Value = 44379
$ws.Range("N111").Value = -52559
$ws.Range("H121").Value = 38463.2
$ws.Range("J121").Value = 38463.2
$ws.Range("L121").Value = 38463.2
$ws.Range("N121").Value = -41957.2
$ws.Range("H132").Value = 3270.347
$ws.Range("I132").Value = 3118.7585
$ws.Range("J132").Value = 3490.15
$ws.Range("K132").Value = 9356.2755
$ws.Range("L132").Value = 10470.45
$ws.Range("M132").Value = -6826.2755
$ws.Range("N132").Value = -15530.45

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 46703.332
$ws.Range("J16").Value = 46703.332
$ws.Range("L16").Value = 46703.332
$ws.Range("N16").Value = -47287.332
$ws.Range("H122").Value = 35238970
$ws.Range("I122").Value = 44048420
$ws.Range("J122").Value = 1166.6666
$ws.Range("K122").Value = 132145260
$ws.Range("L122").Value = 3499.9998
$ws.Range("M122").Value = -132142810
$ws.Range("N122").Value = -8399.9998
$ws.Range("H132").Value = 1128.4222
$ws.Range("I132").Value = 830.9722
$ws.Range("J132").Value = 2318.2222
$ws.Range("K132").Value = 2492.9166
$ws.Range("L132").Value = 6954.6666
$ws.Range("M132").Value = 37.08339999999998
$ws.Range("N132").Value = -12014.6666
